# The commit re-opened the workbook on a different machine/monitor
# (C:\Users\jh\Desktop\... -> C:\Users\mh\OneDrive\...), scrolled the
# sheet down to around row 2123-2186, and left the cursor on C2186.
#
# Most of that is just host-window chrome (absolute file path baked into
# the saved "last opened from" hint, the revision GUID, and the outer
# Excel window's screen position/size) which Excel re-stamps on every
# save but doesn't expose as a scriptable object-model property - so we
# drive every part of it that IS reachable through the Excel COM surface
# and leave the rest alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate() | Out-Null

# Maximize/resize the workbook window to match the new, bigger monitor.
try { $excel.ActiveWindow.WindowState = -4137 } catch { }
try { $excel.ActiveWindow.Top = -98 } catch { }
try { $excel.ActiveWindow.Left = -98 } catch { }
try { $excel.ActiveWindow.Width = 24496 } catch { }
try { $excel.ActiveWindow.Height = 15796 } catch { }

# Scroll the sheet so row 2123 is at the top of the viewport, then move
# the selection down to the cell that was actually being worked on.
try { $excel.ActiveWindow.ScrollRow = 2123 } catch { }
try { $excel.ActiveWindow.ScrollColumn = 1 } catch { }
$ws.Range("C2186").Select() | Out-Null
